$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 44.142857
$ws.Range("I6").Value = 44.833332
$ws.Range("K6").Value = 134.499996
$ws.Range("M6").Value = -22.49999600000001

$ws.Range("H11").Value = 68.85714
$ws.Range("I11").Value = 68.85714
$ws.Range("K11").Value = 68.85714
$ws.Range("M11").Value = 71.14286

$ws.Range("H12").Value = 560.125
$ws.Range("I12").Value = 497
$ws.Range("K12").Value = 497
$ws.Range("M12").Value = -327

$ws.Range("H28").Value = 667628.9
$ws.Range("I28").Value = 1250399.1
$ws.Range("K28").Value = 1250399.1
$ws.Range("M28").Value = -1249914.1

$ws.Range("H53").Value = 252.3125
$ws.Range("I53").Value = 304.81818
$ws.Range("J53").Value = 136.8
$ws.Range("K53").Value = 304.81818
$ws.Range("L53").Value = 136.8
$ws.Range("M53").Value = 332.18182
$ws.Range("N53").Value = -1410.8

$ws.Range("H111").Value = 1804.8
$ws.Range("I111").Value = 1011.6667
$ws.Range("J111").Value = 2994.5
$ws.Range("K111").Value = 3035.0001
$ws.Range("L111").Value = 8983.5
$ws.Range("M111").Value = 31.9998999999998
$ws.Range("N111").Value = -15117.5

$ws.Range("H112").Value = 2265.0625
$ws.Range("J112").Value = 2369.4
$ws.Range("L112").Value = 7108.200000000001
$ws.Range("N112").Value = -9324.200000000001

$ws.Range("H113").Value = 11170
$ws.Range("I113").Value = 8250
$ws.Range("K113").Value = 8250
$ws.Range("M113").Value = -4996

$ws.Range("H137").Value = 611.2
$ws.Range("I137").Value = 636.75
$ws.Range("J137").Value = 509
$ws.Range("K137").Value = 1910.25
$ws.Range("L137").Value = 1527
$ws.Range("M137").Value = 639.75
$ws.Range("N137").Value = -6627

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1572.15
$ws.Range("I63").Value = 1311.9166
$ws.Range("J63").Value = 1962.5
$ws.Range("K63").Value = 1311.9166
$ws.Range("L63").Value = 1962.5
$ws.Range("M63").Value = -625.9166
$ws.Range("N63").Value = -3334.5

$ws.Range("H66").Value = 1572.15
$ws.Range("I66").Value = 1311.9166
$ws.Range("J66").Value = 1962.5
$ws.Range("K66").Value = 6559.583000000001
$ws.Range("L66").Value = 9812.5
$ws.Range("M66").Value = -3127.583000000001
$ws.Range("N66").Value = -16676.5

$ws.Range("H88").Value = 1974.125
$ws.Range("J88").Value = 2005.625
$ws.Range("L88").Value = 2005.625
$ws.Range("N88").Value = -2817.625

$ws.Range("H91").Value = 1974.125
$ws.Range("J91").Value = 2005.625
$ws.Range("L91").Value = 2005.625
$ws.Range("N91").Value = -4813.625

$ws.Range("H97").Value = 577.2353000000001
$ws.Range("I97").Value = 484.2857
$ws.Range("J97").Value = 642.3
$ws.Range("K97").Value = 484.2857
$ws.Range("L97").Value = 642.3
$ws.Range("M97").Value = 11.71429999999998
$ws.Range("N97").Value = -1634.3

$ws.Range("H98").Value = 49833
$ws.Range("J98").Value = 49833
$ws.Range("L98").Value = 49833
$ws.Range("N98").Value = -55823

$ws.Range("H110").Value = 1374.5
$ws.Range("I110").Value = 999.3333
$ws.Range("J110").Value = 2500
$ws.Range("K110").Value = 999.3333
$ws.Range("L110").Value = 2500
$ws.Range("M110").Value = 1045.6667
$ws.Range("N110").Value = -6590

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2888.125
$ws.Range("I134").Value = 2888.125
$ws.Range("K134").Value = 8664.375
$ws.Range("M134").Value = -6129.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3563.5
$ws.Range("I16").Value = 3502
$ws.Range("J16").Value = 3666
$ws.Range("K16").Value = 3502
$ws.Range("L16").Value = 3666
$ws.Range("M16").Value = -3215
$ws.Range("N16").Value = -4240

$ws.Range("H107").Value = 886.4
$ws.Range("I107").Value = 885.5
$ws.Range("K107").Value = 885.5
$ws.Range("M107").Value = 1034.5

$ws.Range("H113").Value = 3563.5
$ws.Range("I113").Value = 3502
$ws.Range("J113").Value = 3666
$ws.Range("K113").Value = 3502
$ws.Range("L113").Value = 3666
$ws.Range("M113").Value = -1332
$ws.Range("N113").Value = -8006

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 574.5
$ws.Range("J34").Value = 800
$ws.Range("L34").Value = 2400
$ws.Range("N34").Value = -2568

$ws.Range("H55").Value = 968.8333
$ws.Range("I55").Value = 966
$ws.Range("K55").Value = 2898
$ws.Range("M55").Value = -2721

$ws.Range("H114").Value = 15649.934
$ws.Range("I114").Value = 225
$ws.Range("J114").Value = 29146.75
$ws.Range("K114").Value = 675
$ws.Range("L114").Value = 87440.25
$ws.Range("M114").Value = 2579
$ws.Range("N114").Value = -93948.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2100
$ws.Range("I113").Value = 1950
$ws.Range("J113").Value = 2300
$ws.Range("K113").Value = 1950
$ws.Range("L113").Value = 2300
$ws.Range("M113").Value = 220
$ws.Range("N113").Value = -6640

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4124.2
$ws.Range("I22").Value = 695
$ws.Range("J22").Value = 6410.3335
$ws.Range("K22").Value = 695
$ws.Range("L22").Value = 6410.3335
$ws.Range("M22").Value = -400
$ws.Range("N22").Value = -7000.3335

$ws.Range("H27").Value = 4124.2
$ws.Range("I27").Value = 695
$ws.Range("J27").Value = 6410.3335
$ws.Range("K27").Value = 695
$ws.Range("L27").Value = 6410.3335
$ws.Range("M27").Value = -588
$ws.Range("N27").Value = -6624.3335

$ws.Range("H46").Value = 1237.2222
$ws.Range("I46").Value = 1088
$ws.Range("J46").Value = 1423.75
$ws.Range("K46").Value = 1088
$ws.Range("L46").Value = 1423.75
$ws.Range("M46").Value = -900
$ws.Range("N46").Value = -1799.75

$ws.Range("H55").Value = 602.3043
$ws.Range("I55").Value = 256.875
$ws.Range("J55").Value = 1391.8572
$ws.Range("K55").Value = 256.875
$ws.Range("L55").Value = 1391.8572
$ws.Range("M55").Value = -83.875
$ws.Range("N55").Value = -1737.8572

$ws.Range("H61").Value = 9207.333000000001
$ws.Range("I61").Value = 9207.333000000001
$ws.Range("K61").Value = 9207.333000000001
$ws.Range("M61").Value = -9005.333000000001

$ws.Range("H113").Value = 9207.333000000001
$ws.Range("I113").Value = 9207.333000000001
$ws.Range("K113").Value = 9207.333000000001
$ws.Range("M113").Value = -7037.333000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 33661.082
$ws.Range("I62").Value = 42641.535
$ws.Range("K62").Value = 42641.535
$ws.Range("M62").Value = -42017.535

$ws.Range("H65").Value = 33661.082
$ws.Range("I65").Value = 42641.535
$ws.Range("K65").Value = 213207.675
$ws.Range("M65").Value = -210087.675

$ws.Range("H81").Value = 5395.857
$ws.Range("I81").Value = 3844.5715
$ws.Range("J81").Value = 6947.143
$ws.Range("K81").Value = 7689.143
$ws.Range("L81").Value = 13894.286
$ws.Range("M81").Value = -6628.143
$ws.Range("N81").Value = -16016.286

$ws.Range("H84").Value = 5395.857
$ws.Range("I84").Value = 3844.5715
$ws.Range("J84").Value = 6947.143
$ws.Range("K84").Value = 38445.715
$ws.Range("L84").Value = 69471.42999999999
$ws.Range("M84").Value = -33141.715
$ws.Range("N84").Value = -80079.42999999999

$ws.Range("H107").Value = 414
$ws.Range("I107").Value = 373.5
$ws.Range("K107").Value = 1120.5
$ws.Range("M107").Value = 799.5

$ws.Range("H122").Value = 4208.727
$ws.Range("I122").Value = 4207.577
$ws.Range("J122").Value = 4213
$ws.Range("K122").Value = 12622.731
$ws.Range("L122").Value = 12639
$ws.Range("M122").Value = -10172.731
$ws.Range("N122").Value = -17539

$ws.Range("H133").Value = 59241.8
$ws.Range("J133").Value = 59241.8
$ws.Range("L133").Value = 59241.8
$ws.Range("N133").Value = -69361.8
